$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$sec.PageSetup.HeaderDistance = 99
Write-Output "done"
